# Replaced with US version 2.1 before dropping in WRI's draft Brazil files
#
# - Renames Sheet1 -> About, BpTPEU -> BpTPEU-large
# - Adds a new BpTPEU-small sheet (clone of BpTPEU-large using 10^3 instead of 10^15)
# - Rewrites the About sheet copy to describe both a "large" and "small" primary
#   energy output unit instead of a single one.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Sheets: rename existing two, clone BpTPEU-large to get BpTPEU-small
#    (this keeps the tab color, column widths and number formats in sync
#    with the "large" sheet, same as the source workbook's own unit sheets).
# ---------------------------------------------------------------------------
$wsAbout = $wb.Worksheets.Item(1)
$wsAbout.Name = "About"

$wsLarge = $wb.Worksheets.Item(2)
$wsLarge.Name = "BpTPEU-large"

$wsLarge.Copy($null, $wsLarge)
$wsSmall = $wb.Worksheets.Item(3)
$wsSmall.Name = "BpTPEU-small"

# ---------------------------------------------------------------------------
# 2. About sheet: clear old layout, write the new content.
# ---------------------------------------------------------------------------
$wsAbout.Cells.Clear()

$wsAbout.Range("A1").Value = "BpTPEU BTU per Large Primary Energy Unit"
$wsAbout.Range("A1").Font.Bold = $true

$wsAbout.Range("A2").Value = "BpTPEU BTU per Small Primary Energy Unit"
$wsAbout.Range("A2").Font.Bold = $true

$wsAbout.Range("A4").Value = "Source:"
$wsAbout.Range("A4").Font.Bold = $true
$wsAbout.Range("B4").Value = "none needed"

$wsAbout.Range("B5").HorizontalAlignment = -4131

$wsAbout.Range("A9").Value = "Notes"
$wsAbout.Range("A9").Font.Bold = $true

$wsAbout.Range("A10").Value = "For the U.S.:"
$wsAbout.Range("A11").Value = "The large primary energy output unit (used in totals graphs) is: quadrillion BTU"
$wsAbout.Range("A12").Value = "The small primary energy output unit (used in energy intensity per unit GDP graphs) is: thousand BTU"

$wsAbout.PageSetup.Orientation = 1

# ---------------------------------------------------------------------------
# 3. BpTPEU-large sheet: relabel B1 (was "One Quadrillion BTU"), left-align it.
#    The A2 "BTU" label and the B2 10^15 formula stay untouched.
# ---------------------------------------------------------------------------
$wsLarge.Range("B1").Value = "large primary energy output unit"
$wsLarge.Range("B1").HorizontalAlignment = -4131

# ---------------------------------------------------------------------------
# 4. BpTPEU-small sheet: same shape as BpTPEU-large but for the "thousand
#    BTU" unit, so the formula becomes 10^3 instead of 10^15.
# ---------------------------------------------------------------------------
$wsSmall.Range("B1").Value = "small primary energy output unit"
$wsSmall.Range("B1").HorizontalAlignment = -4131

$wsSmall.Range("B2").Formula = "=10^3"
$wsSmall.Range("B2").NumberFormat = "General"

# Restore the originally-selected tab (the new sheet becomes active by default).
$wsAbout.Activate()
